$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.281.06'
$ws.Range('E2').Value = '  +0.27%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.608.84'
$ws.Range('E3').Value = '  +0.28%  '
$ws.Range('E4').Value = '  -0.26%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.76'
$ws.Range('E5').Value = '  +0.05%  '
$ws.Range('E6').Value = '  -0.26%  '
$ws.Range('E7').Value = '  +0.48%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.250'
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('E9').Value = '  -0.43%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.29'
$ws.Range('E10').Value = '  +1.11%  '
$ws.Range('E11').Value = '  -0.58%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.832.53'
$ws.Range('E12').Value = '  +0.23%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.613.86'
$ws.Range('E13').Value = '  +0.38%  '
$ws.Range('E14').Value = '  +0.79%  '
$ws.Range('E15').Value = '  +1.04%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '26.286.41'
$ws.Range('E16').Value = '  +0.34%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.59'
$ws.Range('E17').Value = '  +1.65%  '
$ws.Range('E18').Value = '  +0.29%  '
$ws.Range('E19').Value = '  -0.23%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '203.52'
$ws.Range('E20').Value = '  +2.78%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.29'
$ws.Range('E21').Value = '  +1.18%  '
$ws.Range('E22').Value = '  -1.18%  '
$ws.Range('E23').Value = '  -0.17%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.92'
$ws.Range('E24').Value = '  +8.25%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.47'
$ws.Range('E25').Value = '  +1.42%  '
$ws.Range('E26').Value = '  -0.17%  '
$ws.Range('E27').Value = '  -5.84%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.21'
$ws.Range('E28').Value = '  +0.30%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.57'
$ws.Range('E29').Value = '  +1.63%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0488'
$ws.Range('E30').Value = '  +3.73%  '
$ws.Range('E31').Value = '  -0.13%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.19'
$ws.Range('E32').Value = '  +1.95%  '
$ws.Range('E33').Value = '  -2.34%  '
$ws.Range('E34').Value = '  +2.75%  '
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.155.50'
$ws.Range('E36').Value = '  +4.39%  '
$ws.Range('E37').Value = '  +8.55%  '
$ws.Range('E38').Value = '  -0.20%  '
$ws.Range('E39').Value = '  +1.54%  '
$ws.Range('E40').Value = '  +0.48%  '
$ws.Range('E41').Value = '  +0.15%  '
$ws.Range('E42').Value = '  +0.68%  '
$ws.Range('E43').Value = '  +2.76%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.744.73'
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '92.03'
$ws.Range('E45').Value = '  -0.54%  '
$ws.Range('E46').Value = '  -1.24%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '54.26'
$ws.Range('E47').Value = '  +1.06%  '
$ws.Range('E48').Value = '  -0.37%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0₇0985'
$ws.Range('E49').Value = '  -9.21%  '
$ws.Range('E50').Value = '  -0.80%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.00'
$ws.Range('E51').Value = '  -0.46%  '
